# Weekly data refresh: a new price record for this week is inserted at the
# top of the data block (row 112), pushing all the existing records down by
# one row. The oldest record that falls off the bottom of the original range
# (old row 150) ends up as the new last row (151).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 112; this shifts rows 112:150 down to 113:151 and
# extends the sheet dimension to A1:R151 automatically.
$ws.Rows("112:112").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A112").Value = 8
$ws.Range("B112").Value = "Terminal La Palmera de La Serena"
$ws.Range("C112").Value = "Coquimbo"
$ws.Range("D112").Value = 45119
$ws.Range("E112").Value = 4
$ws.Range("F112").Value = 100114007
$ws.Range("G112").Value = "Jengibre"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 200
$ws.Range("K112").Value = 17000
$ws.Range("L112").Value = 18000
$ws.Range("M112").Value = 17500
$ws.Range("N112").Value = "`$/caja 13 kilos"
$ws.Range("O112").Value = "Perú"
$ws.Range("P112").Value = 1346
$ws.Range("Q112").Value = 13
$ws.Range("R112").Value = "Hortaliza"

# Match the date cell's number format to the rest of column D.
$ws.Range("D112").NumberFormat = $ws.Range("D113").NumberFormat
